$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to Text format so numeric-looking values
# (e.g. '1.013') are not reinterpreted as numbers by Excel's smart entry.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.638.18"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.122.13"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "353.45"
$ws.Range("E5").Value = "  +5.55%  "
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "0.5277"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").Value = "0.4531"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "54.13"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "0.09082"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").Value = "1.181"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "24.58"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "2.130.86"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "6.846"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "8.092"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "102.51"
$ws.Range("E16").Value = "  +6.28%  "
$ws.Range("D17").Value = "0.00001176"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "1.013"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "0.06729"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "19.46"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Value = "6.336"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "30.705.90"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "12.81"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").Value = "2.387"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "2.371.73"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "22.48"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "2.575"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "165.28"
$ws.Range("E29").Value = "  +1.10%  "
$ws.Range("D30").Value = "137.08"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "1.197"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "1.654"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "6.371"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "4.028"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "6.007"
$ws.Range("E36").Value = "  +6.25%  "
$ws.Range("D37").Value = "10.32"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "0.02658"
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("D39").Value = "0.06879"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D40").Value = "0.2320"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").Value = "12.55"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("D42").Value = "0.6916"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "1.280"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("D44").Value = "14.81"
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.6473"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "2.332"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "3.778"
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("D48").Value = "0.00000000366"
$ws.Range("E48").Value = "  +8.59%  "
$ws.Range("D49").Value = "1.257"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "82.89"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "0.07312"
$ws.Range("E51").Value = "  +2.29%  "
